$wb = $excel.ActiveWorkbook

# --- Parametric survival model estimates (per-arm) ---
# weibull
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.95410791709162
$ws.Range("C2").Value = 0.173133254677496
$ws.Range("B3").Value = 0.0556593266224006
$ws.Range("C3").Value = 0.0926145582214106

# lognormal
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.40219799434963
$ws.Range("C2").Value = 0.234519749464493
$ws.Range("B3").Value = -1.00901173307717
$ws.Range("C3").Value = 0.102287905033608

# llogis
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.35918673923445
$ws.Range("C2").Value = 0.116437894083279
$ws.Range("B3").Value = 1.75657302549737
$ws.Range("C3").Value = 0.210949146448444

# gompertz
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.60453622180521
$ws.Range("C2").Value = 0.140775675694235
$ws.Range("B3").Value = -0.0179874448132232
$ws.Range("C3").Value = 0.00981579720341566

# --- Covariance matrices (per model) ---
# weibull cov
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0299751238752225
$ws.Range("B2").Value = -0.0112536040524681
$ws.Range("A3").Value = -0.0112536040524681
$ws.Range("B3").Value = 0.00857745639454706

# lognormal cov
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0549995128888884
$ws.Range("B2").Value = -0.0211596960556413
$ws.Range("A3").Value = -0.0211596960556413
$ws.Range("B3").Value = 0.0104628155161643

# llogis cov
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.013557783178549
$ws.Range("B2").Value = 0.00663085164219347
$ws.Range("A3").Value = 0.00663085164219347
$ws.Range("B3").Value = 0.0444995423873269

# gompertz cov
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0198177908671686
$ws.Range("B2").Value = -0.000793480139508984
$ws.Range("A3").Value = -0.000793480139508984
$ws.Range("B3").Value = 0.0000963498747385826
